$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The script that produces this report re-ran for 2025-05-22: the day that
# used to be the last row (115, with "NA" in the "Numero de page" column)
# is no longer the last one, so its "NA" moves down to the freshly appended
# row 116, and row 115's "Numero de page" cell is cleared out.
$ws.Range("C115").Value = ""

# Force column A to stay plain text so the date string isn't reinterpreted
# as a date serial number, matching the rest of the "Date" column.
$ws.Range("A116").NumberFormat = "@"
$ws.Range("A116").Value = "2025-05-22"
$ws.Range("A116").ClearFormats()

$ws.Range("B116").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C116").Value = "NA"
$ws.Range("D116").Value = 1
